# Natmi following Dr Hou advice
# Recompute the Thbs1 -> Tnfrsf11b LR-pairs table: instead of a single
# "Target cluster" (FAPs) for every "Sending cluster", each sending
# cluster now gets two rows - one for Target cluster FAPs and one for
# Target cluster sCs - with refreshed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="ECs";  D="FAPs"; E=3; F=1; G=122.253015;         H=366.759045;         I=0.1988639364328829; J=0.1988639364328829; K=3; L=1;                  M=3.776574666666666;  N=11.329724;  O=0.9855052394405499;  P=0.9855052394405499;  Q=461.69763937262;     R=4155.27875435358;   S=0.1959814512903786;  T=0.1959814512903785 },
    @{ Row=3; A="ECs";  D="sCs";  E=3; F=1; G=122.253015;         H=366.759045;         I=0.1988639364328829; J=0.1988639364328829; K=1; L=0.3333333333333333; M=0.05554566666666667; N=0.166637;   O=0.01449476055945007; P=0.01449476055945008; Q=6.790625220185;      R=61.115626981665;    S=0.002882485142504337; T=0.002882485142504337 },
    @{ Row=4; A="FAPs"; D="FAPs"; E=3; F=1; G=132.5447616666667;  H=397.634285;         I=0.2156050961899926; J=0.2156050961899926; K=3; L=1;                  M=3.776574666666666;  N=11.329724;  O=0.9855052394405499;  P=0.9855052394405499;  Q=500.5651891097044;   R=4505.086701987339;  S=0.2124799519453215;  T=0.2124799519453215 },
    @{ Row=5; A="FAPs"; D="sCs";  E=3; F=1; G=132.5447616666667;  H=397.634285;         I=0.2156050961899926; J=0.2156050961899926; K=1; L=0.3333333333333333; M=0.05554566666666667; N=0.166637;   O=0.01449476055945007; P=0.01449476055945008; Q=7.362287149949444;   R=66.260584349545;    S=0.003125144244671144; T=0.003125144244671144 },
    @{ Row=6; A="M2";   D="FAPs"; E=3; F=1; G=320.0894206666666;  H=960.2682619999999;  I=0.5206762565675317; J=0.5206762565675317; K=3; L=1;                  M=3.776574666666666;  N=11.329724;  O=0.9855052394405499;  P=0.9855052394405499;  Q=1208.841597157743;   R=10879.57437441969;  S=0.5131291788995945;  T=0.5131291788995945 },
    @{ Row=7; A="M2";   D="sCs";  E=3; F=1; G=320.0894206666666;  H=960.2682619999999;  I=0.5206762565675317; J=0.5206762565675317; K=1; L=0.3333333333333333; M=0.05554566666666667; N=0.166637;   O=0.01449476055945007; P=0.01449476055945008; Q=17.77958026387711;   R=160.016222374894;   S=0.007547077667937166; T=0.007547077667937167 },
    @{ Row=8; A="sCs";  D="FAPs"; E=3; F=1; G=39.86989333333333;  H=119.60968;          I=0.06485471080959287; J=0.06485471080959287; K=3; L=1;                 M=3.776574666666666;  N=11.329724;  O=0.9855052394405499;  P=0.9855052394405499;  Q=150.5716291253689;   R=1355.14466212832;   S=0.06391465730525545; T=0.06391465730525545 },
    @{ Row=9; A="sCs";  D="sCs";  E=3; F=1; G=39.86989333333333;  H=119.60968;          I=0.06485471080959287; J=0.06485471080959287; K=1; L=0.3333333333333333; M=0.05554566666666667; N=0.166637; O=0.01449476055945007; P=0.01449476055945008; Q=2.214599805128889;   R=19.93139824616;     S=0.0009400535043374271; T=0.0009400535043374272 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = $r.A
    $ws.Range("B$i").Value = "Thbs1"
    $ws.Range("C$i").Value = "Tnfrsf11b"
    $ws.Range("D$i").Value = $r.D
    $ws.Range("E$i").Value = $r.E
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = $r.G
    $ws.Range("H$i").Value = $r.H
    $ws.Range("I$i").Value = $r.I
    $ws.Range("J$i").Value = $r.J
    $ws.Range("K$i").Value = $r.K
    $ws.Range("L$i").Value = $r.L
    $ws.Range("M$i").Value = $r.M
    $ws.Range("N$i").Value = $r.N
    $ws.Range("O$i").Value = $r.O
    $ws.Range("P$i").Value = $r.P
    $ws.Range("Q$i").Value = $r.Q
    $ws.Range("R$i").Value = $r.R
    $ws.Range("S$i").Value = $r.S
    $ws.Range("T$i").Value = $r.T
}
